$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demo data")

# --- Text cell updates (shared string changes) ---
$ws.Range("AO8").Value = "Principal Agent:  This them that"
$ws.Range("AL9").Value = "Principal Agent:  Fred’s Steel Manufacturing"

# --- Latitude (column H) value tweaks ---
$ws.Range("H2").Value = -30.13029
$ws.Range("H3").Value = -29.86148
$ws.Range("H4").Value = -30.1057
$ws.Range("H5").Value = -29.54888
$ws.Range("H7").Value = -32.23235
$ws.Range("H8").Value = -31.0595994
$ws.Range("H9").Value = -31.5727543
$ws.Range("H10").Value = -24.8472
$ws.Range("H11").Value = -32.66633
$ws.Range("H12").Value = -26.74655

# --- View/selection state: scroll back to top-left and move selection to H13 ---
$ws.Activate()
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H13").Select()
